$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Team" column (foreign-key-style column pointing every row at the "UVA" team)
# Header cell: same bold/centered header look as columns B1:I1, but with a
# left+right thin box border instead of the full box used on the other headers.
$ws.Range("J1").Value = "Team"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("J1").VerticalAlignment = -4160     # xlTop
$ws.Range("J1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft, xlContinuous
$ws.Range("J1").Borders.Item(10).LineStyle = 1  # xlEdgeRight, xlContinuous

# Every roster row belongs to the same team
$ws.Range("J2:J105").Value = "UVA"

# Keep the sheet's selection in sync with the now-wider used range
$ws.Range("A1:J105").Select() | Out-Null
